# Journal de bord - Jungo
# Add a new journal entry (row 16) to the "Feuil1" worksheet:
#   Date  = 2020-03-27 (serial 43917)
#   Temps = 1.5
#   Description = "Entrevue avec le chef de projet, il m'a donné du code
#                  pour que je puisse avancer"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Write the new entry's values --------------------------------------
$ws.Range("A16").Value = 43917
$ws.Range("B16").Value = 1.5
$ws.Range("C16").Value = "Entrevue avec le chef de projet, il m'a donné du code pour que je puisse avancer"

# --- Match the formatting of the row above it (row 14 is a similarly
#     wrapped, two-line entry) by copying its cell formats in one shot,
#     instead of toggling individual properties (which would create
#     spurious intermediate cell styles). -------------------------------
$ws.Range("A14:C14").Copy()
$ws.Range("A16:C16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 14 (like the other wrapped, multi-line entries) uses a taller row
# height than the sheet default.
$ws.Rows.Item(16).RowHeight = $ws.Rows.Item(14).RowHeight

# --- Move the active selection to the newly added cell -----------------
$ws.Range("C16").Select()
